$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.072.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.08%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.124.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.52%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'593.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.57%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'136.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.02%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.116.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.60%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.07%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.78%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.80%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -3.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'34.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.91%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.635.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.59%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.37%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.065.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.02%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.121.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.75%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'472.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.96%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.696"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.44%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'86.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.85%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'12.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.97%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.00%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'26.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "'  -6.93%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.26%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.50%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -3.20%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'52.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.81%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -7.65%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'422.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0386"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.42%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.95%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -10.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.899.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.113"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.264"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.56%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.41%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'25.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.18%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.57%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -5.74%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'120.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.25%  "
$ws.Range("E51").Style = "Normal"
